# chapter 5 initial draft
#
# The authored change inserts a brand-new first column ("Cd3" / 0.02) in
# front of the existing Sheet1 data table, pushing the other seven
# columns (Dynamic Pressure, mFuel2, m2, Cd, Isp, T3, m3) one place to
# the right (A:G -> B:H), and widens the bar chart's series references
# to match the new A:H extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing headers
# (row 1) and values (row 2) from A:G to B:H, matching the sheet1.xml
# diff (dimension A1:G2 -> A1:H2, spans 1:7 -> 1:8, each column letter
# shifted up by one).
$ws.Columns("A").Insert()

# Populate the freshly inserted column with the new header/value pair.
$ws.Range("A1").Value = "Cd3"
$ws.Range("A2").Value = 0.02

# The bar chart's single series still points at the old A:G extent;
# repoint both the category (header) and value references at the new
# A:H extent so the chart keeps tracking the whole table.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = '=SERIES(,Sheet1!$A$1:$H$1,Sheet1!$A$2:$H$2,1)'

# Match the author's final cursor position captured in the workbook.
$ws.Range("L4").Select()
